$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -21.756
$ws.Range("D4").Value = -8.132999999999999

$ws.Range("D5").Value = -8.581999999999999

$ws.Range("A6").Value = -21.108

$ws.Range("A7").Value = -21.047

$ws.Range("D8").Value = -8.203999999999999

$ws.Range("A16").Value = -20.714
$ws.Range("D16").Value = -8.463999999999999

$ws.Range("A20").Value = -21.98

$ws.Range("D22").Value = -8.16
